$d = $word.ActiveDocument

# --- Paragraph 3: "In the current project the 2 adapters ..." ---
# Only proofing-mark/run-merge cleanup; text content unchanged.
$p3 = $d.Paragraphs(3)
$start = $p3.Range.Start
$end = $p3.Range.End
$d.Range($start, $end).Text = "In the current project the 2 adapters were made using a similar logic : FormattedPage, FormattedAlbum. Those components compose a Page object and an Album object respectively. They only alter the original classes' ToString method so it will be more human readable and suit the ListBox UI. Therefore it is considered a use of the adapter design pattern."

# --- Paragraph 7: "In the current project it will be used in order to create the 'UserAverageableDetails' ..." ---
$p7 = $d.Paragraphs(7)
$start = $p7.Range.Start
$end = $p7.Range.End
$d.Range($start, $end).Text = "In the current project it will be used in order to create the 'UserAverageableDetails' objects (also the builder component) using the 'Averageizer' object.  This usage is correct since 'averageable details' is not a coherent definition and could be easily expanded, diminished or altered in the future. This may suggest that in the future that the class structure will be changed. For example, a user's favorite musical genre can be considered as 'avereageable': "

# --- Paragraph 8: "MusicAverage(SoulMusic , Jazz) = Funk." ---
# Keep the special "S" run (with cs font hint) untouched; merge the runs before and after it.
$p8 = $d.Paragraphs(8)
$p8start = $p8.Range.Start
$p8full = $p8.Range.Text
$sIdx = $p8full.IndexOf("S")
$beforeStart = $p8start
$beforeEnd = $p8start + $sIdx
$d.Range($beforeStart, $beforeEnd).Text = "MusicAverage("
# Recompute paragraph 8 after the first edit (length may have changed, though here it is identical length)
$p8 = $d.Paragraphs(8)
$p8full = $p8.Range.Text
$sIdx = $p8full.IndexOf("S")
$afterStart = $p8.Range.Start + $sIdx + 1
$afterEnd = $p8.Range.End
$d.Range($afterStart, $afterEnd).Text = "oulMusic , Jazz) = Funk."

# --- Paragraph 10: "Also, the composer may be altered, for instance the 'ClosestCity' ..." ---
$p10 = $d.Paragraphs(10)
$start = $p10.Range.Start
$end = $p10.Range.End
$d.Range($start, $end).Text = "Also, the composer may be altered, for instance the 'ClosestCity' property may be calculated in a different way rather than a simple aerial average."

# --- Paragraph 13: "In the current project the cities data-base was created ..." ---
# Only the first run changes (an extra sentence is inserted); the remainder of the
# paragraph (with unique details / Thus a single source... / Facade sentence) stays.
$p13 = $d.Paragraphs(13)
$full = $p13.Range.Text
$oldFirst = "In the current project the cities data-base was created in a separate file as a singleton. This creation method is highly important since the cities supposed to have a unique instance "
$idx = $full.IndexOf($oldFirst)
$start = $p13.Range.Start + $idx
$end = $start + $oldFirst.Length
$d.Range($start, $end).Text = "In the current project the cities data-base was created in a separate file as a singleton- the 'CitiesDataBase' class. This creation method is highly important since the cities supposed to have a unique instance "
